$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last column (BA, old "Mean" column with the 51-run mean
# that included "Run 50"). This shifts nothing else - BA was the last
# column - and drops the dimension from A1:BA14 to A1:AZ14.
$ws.Columns("BA:BA").Delete()

# Column A header: "Gen" -> "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# Column A body: generation counts -> fraction-of-budget values
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 0.001
$ws.Range("A4").Value = 0.01
$ws.Range("A5").Value = 0.1
$ws.Range("A6").Value = 0.2
$ws.Range("A7").Value = 0.3
$ws.Range("A8").Value = 0.4
$ws.Range("A9").Value = 0.5
$ws.Range("A10").Value = 0.6
$ws.Range("A11").Value = 0.7
$ws.Range("A12").Value = 0.8
$ws.Range("A13").Value = 0.9
$ws.Range("A14").Value = 1

# Column AZ (now the last column, previously "Run 50") becomes the new
# "Mean" column, recomputed over the 50 remaining runs (B:AY).
$ws.Range("AZ1").Value = "Mean"
$ws.Range("AZ2").Value = 16122753739.30056
$ws.Range("AZ3").Value = 11503974516.30149
$ws.Range("AZ4").Value = 922269946.0993564
$ws.Range("AZ5").Value = 4063.24595453
$ws.Range("AZ6").Value = 1123.88527143
$ws.Range("AZ7").Value = 591.06965975
$ws.Range("AZ8").Value = 467.18684602
$ws.Range("AZ9").Value = 435.42107295
$ws.Range("AZ10").Value = 413.62585038
$ws.Range("AZ11").Value = 406.62464817
$ws.Range("AZ12").Value = 403.24245822
$ws.Range("AZ13").Value = 401.89840072
$ws.Range("AZ14").Value = 401.3468612
